$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: mirror "Yes" (same as B9) into C9 and D9 ---
$ws.Range("C9").Value = "Yes"
$ws.Range("D9").Value = "Yes"

# --- New IPDC benefits/VfM rows (284-290) ---

# Row 284: label only (Arial-styled label cell)
$ws.Range("A284").Value = "Project stage"
$ws.Range("A284").Font.Name = "Arial"
$ws.Range("A284").Font.Size = 10

# Row 285: label (Arial) + three numeric values
$ws.Range("A285").Value = "Initial Benefits Cost Ratio (BCR)"
$ws.Range("A285").Font.Name = "Arial"
$ws.Range("A285").Font.Size = 10
$ws.Range("B285").Value = 609
$ws.Range("C285").Value = 2289
$ws.Range("D285").Value = 82798

# Row 286: label only (Arial)
$ws.Range("A286").Value = "Adjusted Benefits Cost Ratio (BCR)"
$ws.Range("A286").Font.Name = "Arial"
$ws.Range("A286").Font.Size = 10

# Row 287: label only (Arial)
$ws.Range("A287").Value = "VfM Category single entry"
$ws.Range("A287").Font.Name = "Arial"
$ws.Range("A287").Font.Size = 10

# Row 288: label only (default style, no explicit Arial font)
$ws.Range("A288").Value = "VfM Category lower range"

# Row 289: label only (Arial)
$ws.Range("A289").Value = "VfM Category upper range"
$ws.Range("A289").Font.Name = "Arial"
$ws.Range("A289").Font.Size = 10

# Row 290: label (default style) + Red/Green values
$ws.Range("A290").Value = "SRO Benefits RAG"
$ws.Range("B290").Value = "Red"
$ws.Range("C290").Value = "Green"
$ws.Range("D290").Value = "Green"

# --- View state: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 279
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D290").Select()
